# Tumor.xlsx update
#
# Sheet "B" (the second worksheet) stored its "marked" cells as a shared
# string ("M") and its "unmarked" cells as the number 0. This pass flips
# every cell in B2:J10 on sheet "B": cells that used to read "M" become the
# number 0, and cells that used to read 0 become the number 1 - i.e. the
# whole grid is inverted and re-written with plain numbers so the
# now-unused "M" shared string disappears from the workbook. The cursor
# position on both sheets is also updated to reflect where the user left
# off working.

$wb = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item(1)
$wsB = $wb.Worksheets.Item(2)

# New (inverted) B2:J10 grid for sheet "B", written as plain numbers.
$data = @(
  @(1, 0, 1, 1, 1, 0, 1, 1, 1),
  @(1, 0, 0, 0, 1, 0, 1, 0, 1),
  @(1, 1, 1, 0, 1, 0, 0, 1, 1),
  @(1, 1, 0, 0, 0, 0, 1, 1, 1),
  @(1, 1, 1, 0, 0, 1, 1, 1, 1),
  @(1, 1, 1, 1, 0, 0, 0, 1, 1),
  @(1, 1, 1, 1, 0, 1, 0, 0, 1),
  @(1, 1, 1, 1, 1, 1, 1, 1, 1),
  @(1, 1, 1, 1, 1, 1, 1, 1, 1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $wsB.Cells.Item($r + 2, $c + 2).Value = $row[$c]
  }
}

# Restore the cursor positions that were left selected on each sheet.
$wsA.Range("N10").Select()
$wsB.Activate()
$wsB.Range("G13").Select()
